$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule": update Cost ($) and Unit Cost ($/ML) for rows 4-6 ---
$wsSchedule = $wb.Worksheets.Item("Schedule")

$wsSchedule.Range("E4").Value = 660.1421482500001
$wsSchedule.Range("F4").Value = 43.66019499007938

$wsSchedule.Range("E5").Value = 564.6798300000002
$wsSchedule.Range("F5").Value = 21.34088548752835

$wsSchedule.Range("E6").Value = 408.44997375
$wsSchedule.Range("F6").Value = 27.01388715277778

# --- Sheet "Detailed": update Price and Type for rows 38-97 ---
$wsDetailed = $wb.Worksheets.Item("Detailed")

$wsDetailed.Range("B38").Value = 100.01
$wsDetailed.Range("B39").Value = 108.89
$wsDetailed.Range("B40").Value = 109.12055
$wsDetailed.Range("C40").Value = "historical"
$wsDetailed.Range("B41").Value = 112.06343
$wsDetailed.Range("C41").Value = "historical"
$wsDetailed.Range("B42").Value = 109.91454
$wsDetailed.Range("C42").Value = "historical"
$wsDetailed.Range("B43").Value = 105
$wsDetailed.Range("C43").Value = "historical"
$wsDetailed.Range("B44").Value = 98.92904
$wsDetailed.Range("C44").Value = "historical"
$wsDetailed.Range("B45").Value = 84.79000000000001
$wsDetailed.Range("C45").Value = "historical"
$wsDetailed.Range("B46").Value = 70.57004999999999
$wsDetailed.Range("C46").Value = "historical"
$wsDetailed.Range("B47").Value = 99.56278
$wsDetailed.Range("C47").Value = "historical"
$wsDetailed.Range("B48").Value = 84.79000000000001
$wsDetailed.Range("C48").Value = "historical"
$wsDetailed.Range("B49").Value = 85.65000000000001
$wsDetailed.Range("C49").Value = "historical"
$wsDetailed.Range("B50").Value = 79.57692
$wsDetailed.Range("B51").Value = 73.20008
$wsDetailed.Range("B52").Value = 73.20009
$wsDetailed.Range("B53").Value = 73.20009
$wsDetailed.Range("B54").Value = 73.20009
$wsDetailed.Range("B55").Value = 78
$wsDetailed.Range("B56").Value = 78
$wsDetailed.Range("B57").Value = 78
$wsDetailed.Range("B58").Value = 78
$wsDetailed.Range("B59").Value = 66.62696
$wsDetailed.Range("B60").Value = 79.39904
$wsDetailed.Range("B61").Value = 79.95014
$wsDetailed.Range("B62").Value = 96.03658
$wsDetailed.Range("B63").Value = 105.79
$wsDetailed.Range("B64").Value = 94.22657
$wsDetailed.Range("B65").Value = 57.06009
$wsDetailed.Range("B66").Value = 55.92455
$wsDetailed.Range("B67").Value = 46.17584
$wsDetailed.Range("B68").Value = 37.89016
$wsDetailed.Range("B69").Value = 38.73791
$wsDetailed.Range("B70").Value = 37.89023
$wsDetailed.Range("B71").Value = 37.89022
$wsDetailed.Range("B72").Value = 37.89021
$wsDetailed.Range("B73").Value = 40.24846
$wsDetailed.Range("B74").Value = 37.89021
$wsDetailed.Range("B75").Value = 37.89038
$wsDetailed.Range("B76").Value = 37.89018
$wsDetailed.Range("B77").Value = 37.89017
$wsDetailed.Range("B78").Value = 37.89019
$wsDetailed.Range("B79").Value = 64.89
$wsDetailed.Range("B80").Value = 84.79000000000001
$wsDetailed.Range("B81").Value = 80.50707
$wsDetailed.Range("B82").Value = 91.21869
$wsDetailed.Range("B83").Value = 78
$wsDetailed.Range("B84").Value = 84.79000000000001
$wsDetailed.Range("B85").Value = 78
$wsDetailed.Range("B86").Value = 69.44401000000001
$wsDetailed.Range("B87").Value = 70.36225
$wsDetailed.Range("B88").Value = 71.02005
$wsDetailed.Range("B89").Value = 70.35384000000001
$wsDetailed.Range("B90").Value = 69.03394
$wsDetailed.Range("B91").Value = 57.31
$wsDetailed.Range("B92").Value = 57.06
$wsDetailed.Range("B95").Value = 56.98
$wsDetailed.Range("B97").Value = 45.77911
